# Update the crypto price/volume table with the latest scraped values.
# Column D ("Price") values are written as plain text (the sheet always
# stored these as text, e.g. "70.929.74" or "1.00"), so we force the
# cell's number format to Text ("@") before assigning and restore the
# default "Normal" style afterwards so Excel does not silently convert
# the numeric-looking strings into real numbers (which would drop
# formatting such as trailing zeros or thousands separators).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "70.929.74"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +5.65%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.646.02"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +5.33%  "
$ws.Range("E4").Value = "  -0.09%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "594.68"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +2.19%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "195.60"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +2.46%  "
$ws.Range("E7").Value = "  +2.44%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "3.640.81"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +5.39%  "
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("E10").Value = "  +7.26%  "
$ws.Range("E11").Value = "  +4.50%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "58.30"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +1.29%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.0000297"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +6.45%  "
$ws.Range("E14").Value = "  +4.92%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "4.231.16"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +5.46%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "20.39"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +7.53%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "3.648.77"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +5.54%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "70.928.07"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +5.60%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "12.78"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +5.18%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.122"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +2.18%  "
$ws.Range("E21").Value = "  +3.55%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "489.57"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +1.50%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "19.44"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +15.15%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "5.26"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -1.50%  "
$ws.Range("E25").Value = "  +2.77%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "91.51"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +1.45%  "
$ws.Range("E27").Value = "  +5.55%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "11.44"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +3.83%  "
$ws.Range("E29").Value = "  +5.82%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "32.83"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +4.48%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "7.89"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +5.14%  "
$ws.Range("E32").Value = "  +9.78%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "12.29"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +3.72%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "616.96"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +2.14%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "66.46"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +3.27%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "40.31"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +7.10%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.0₃0834"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +9.88%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.413"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +4.94%  "
$ws.Range("E39").Value = "  +0.32%  "
$ws.Range("E40").Value = "  -0.10%  "
$ws.Range("E41").Value = "  +2.26%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "3.329.25"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +3.28%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "3.23"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +13.39%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "3.17"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +6.96%  "
$ws.Range("E45").Value = "  +9.05%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.0459"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +5.58%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "9.64"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +10.23%  "
$ws.Range("E48").Value = "  +2.71%  "
$ws.Range("E49").Value = "  +2.95%  "
$ws.Range("B50").Value = "LidoDAOToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "3.24"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.82%  "
$ws.Range("B51").Value = "FirstDigitalUSD"
$ws.Range("C51").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.07%  "
